# Renames of jcoin:final_variable_name values + column width tweaks
# for table-schema-baseline (col L) and table-schema-time-points (col K)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("table-schema-baseline")
$ws1Updates = @{
    4 = "state_of_site_enrollment"
    6 = "birth_date"
    7 = "age"
    8 = "sex_at_birth"
    9 = "gender_id"
    10 = "gender_id_condensed"
    13 = "race_AIAN"
    14 = "race_hawaiian_OPI"
    17 = "race_AI_tribe"
    18 = "race_other_specified"
    19 = "hispanic_latino"
    20 = "sex_orient_category"
    21 = "sex_orient_other"
    22 = "ever_pregnant"
    24 = "living_as_married"
    25 = "educ_category"
    26 = "educ_highest_grade"
    27 = "educ_other_specified"
    28 = "intv_while_incarc"
    29 = "days_incarcerated_interval"
    30 = "ever_rx_moud"
    31 = "months_daily_bup"
    32 = "months_sublocade"
    33 = "months_weekly_brixadi"
    34 = "months_monthly_brixadi"
    35 = "months_probuphine_implant"
    36 = "months_daily_ntx"
    37 = "months_monthly_vivitrol"
    38 = "months_methadone"
}
foreach ($row in $ws1Updates.Keys) {
    $ws1.Range("L$row").Value = $ws1Updates[$row]
}
$ws1.Columns.Item(12).ColumnWidth = 25.166666666666668

$ws2 = $wb.Worksheets.Item("table-schema-time-points")
$ws2Updates = @{
    4 = "shifted_visit_date"
    7 = "last_getting_using_drugs"
    8 = "last_social_problems"
    9 = "last_work_life_disruption"
    10 = "last_withdrawal"
    11 = "last_used_opioids"
    12 = "last_opioid_overdose"
    13 = "last_went_moud"
    14 = "times_opioid_overdose"
    15 = "times_naloxone"
    16 = "who_gave_naloxone"
    17 = "drugs_before_overdose"
    18 = "times_ems_after_overdose"
    19 = "times_er_after_overdose"
    20 = "times_hospital_after_overdose"
    21 = "tx_refer_after_overdose"
    22 = "days_drug_use"
    23 = "days_alcohol_use"
    24 = "days_alcohol_binge"
    25 = "days_medical_marijuana"
    26 = "days_other_marijuana"
    27 = "days_heroin"
    28 = "days_fentanyl"
    29 = "days_non_rx_methadone"
    30 = "days_non_rx_suboxone"
    31 = "days_other_opioids"
    32 = "days_cocaine"
    33 = "days_amphetamines"
    34 = "days_benzo"
    35 = "days_other_drugs"
    36 = "other_drug_describe"
    37 = "days_confined_no_use"
    38 = "days_illegal_activity"
    39 = "times_drug_possession"
    40 = "times_public_drunk"
    41 = "times_dui"
    42 = "times_drug_dealing"
    43 = "times_vandalism"
    44 = "times_stolen_goods"
    45 = "times_forgery"
    46 = "times_shoplift"
    47 = "times_theft"
    48 = "times_burglary"
    49 = "times_auto_theft"
    50 = "times_carjacking"
    51 = "times_assault"
    52 = "times_robbery"
    53 = "times_aggravated_assault"
    54 = "times_rape"
    55 = "times_homicide"
    56 = "times_arson"
    57 = "times_prostitution"
    58 = "times_other_unlawful"
    59 = "number_charged_arrests"
    60 = "number_arrest_possession"
    61 = "number_arrest_drunk"
    62 = "number_arrest_dui"
    63 = "number_arrest_drug_dealing"
    64 = "number_arrest_vandalism"
    65 = "number_arrest_stolen_goods"
    66 = "number_arrest_forgery"
    67 = "number_arrest_shoplifting"
    68 = "number_arrest_theft"
    69 = "number_arrest_burglary"
    70 = "number_arrest_auto_theft"
    71 = "number_arrest_carjacking"
    72 = "number_arrest_assault"
    73 = "number_arrest_robbery"
    74 = "number_arrest_aggravated"
    75 = "number_arrest_rape"
    76 = "number_arrest_homicide"
    77 = "number_arrest_arson"
    78 = "number_arrest_prostitution"
    79 = "number_arrest_other"
    80 = "days_electronic_monitoring"
    81 = "days_house_arrest"
    82 = "days_jail"
    83 = "days_prison"
    84 = "currently_incarcerated"
    85 = "length_current_incarceration"
    86 = "days_parole"
    87 = "days_probation"
    88 = "days_other_supervision"
    89 = "days_met_probation_officer"
    90 = "days_trouble_probation_officer"
    91 = "number_lifetime_arrests"
    92 = "age_first_arrest"
    93 = "years_lifetime_incarceration"
    94 = "months_lifetime_incarceration"
    95 = "times_guilty_sentenced"
    96 = "age_first_convicted"
    97 = "times_er_visits"
    98 = "nights_hospital_detox"
    99 = "nights_hospitalized"
    100 = "nights_residential_detox"
    101 = "nights_residential_sud_tx"
    102 = "nights_residential_mh"
    103 = "nights_physical_rehab"
    104 = "visits_primary_care"
    105 = "pc_reason_alc_drugs"
    106 = "pc_reason_mh"
    107 = "pc_reason_physical"
    108 = "pc_reason_other"
    109 = "pc_reason_other_specify"
    110 = "op_tx_days"
    111 = "op_tx_days_visited"
    112 = "op_tx_days_online"
    113 = "op_tx_days_doctor"
    114 = "op_tx_days_therapy"
    115 = "op_tx_days_moud_only"
    116 = "psych_visits"
    117 = "psych_visits_in_person"
    118 = "psych_visits_online"
    119 = "counselor_visits"
    120 = "counselor_visits_in_person"
    121 = "counselor_visits_online"
    122 = "healthcare_expense"
    123 = "received_sud_tx"
    124 = "sud_tx_organized"
    125 = "sud_tx_satisfied"
    126 = "sud_tx_efficient"
    127 = "sud_tx_personal"
    128 = "sud_tx_moud"
    129 = "oud_tx_not_appropriate"
    130 = "prefer_oud_meds"
    131 = "prefer_oud_detox"
    132 = "prefer_oud_op"
    133 = "prefer_oud_iop"
    134 = "prefer_oud_residential"
    135 = "prefer_oud_other_tx"
    136 = "prefer_oud_other_specify"
    137 = "prefer_oud_no_treatment"
    138 = "prefer_oud_dk"
    139 = "prefer_moud_type"
    140 = "prefer_which_bup"
    141 = "prefer_which_naloxone"
    142 = "household_people"
    143 = "household_under_18"
    144 = "household_income"
    145 = "household_income_legal"
    146 = "household_public_assist_any"
    147 = "household_public_assist_dollars"
    148 = "household_nonemploy_any"
    149 = "household_nonemploy_dollars"
    150 = "household_income_illegal_any"
    151 = "household_income_illegal_dollars"
    152 = "current_work_school"
    153 = "work_days"
    154 = "work_typical_days_week"
    155 = "work_hours_week"
    156 = "work_hourly_wage"
    157 = "work_health_insurance"
    158 = "work_paid_time_off"
    159 = "work_pension"
    160 = "work_retirement_plan"
    161 = "work_occupation"
    162 = "have_health_insurance"
    163 = "insured_private"
    164 = "insured_medicare"
    165 = "insured_medigap"
    166 = "insured_medicaid"
    167 = "insured_chip"
    168 = "insured_military"
    169 = "insured_indian"
    170 = "insured_state"
    171 = "insured_other_gov"
    172 = "insured_single_service"
    173 = "insured_dont_know"
    174 = "days_uninsured"
    175 = "days_self_help"
    176 = "days_drugfree_activity"
    177 = "days_homeless"
    178 = "days_shelter"
    179 = "days_household_alcohol"
    180 = "days_household_drug"
    181 = "days_activities_substance_use"
    182 = "days_family_trouble"
    183 = "days_arguments"
    184 = "days_abused"
    185 = "narcan_received"
    186 = "narcan_used"
    187 = "narcan_refilled"
    188 = "difficulty_chores"
    189 = "difficulty_stairs"
    190 = "difficulty_walking"
    191 = "difficulty_traveling"
    192 = "trouble_with_leisure"
    193 = "trouble_with_family"
    194 = "trouble_with_work"
    195 = "trouble_with_activities"
    196 = "past_week_fearful"
    197 = "past_week_anxiety"
    198 = "past_week_worried"
    199 = "past_week_uneasy"
    200 = "past_week_worthless"
    201 = "past_week_helpless"
    202 = "past_week_depressed"
    203 = "past_week_hopeless"
    204 = "past_week_fatigued"
    205 = "past_week_tired"
    206 = "past_week_rundown"
    207 = "fatigue_level"
    208 = "sleep_quality"
    209 = "sleep_refreshing"
    210 = "sleep_problems"
    211 = "sleep_difficulty"
    212 = "can_concentrate"
    213 = "can_remember"
    214 = "pain_daily_activity"
    215 = "pain_work_around_house"
    216 = "pain_social_activity"
    217 = "pain_household_chores"
    218 = "pain_intensity"
    219 = "last_time_multiple_partners"
    220 = "last_time_unprotected_sex"
    221 = "last_time_sex_intoxicated"
    222 = "last_time_injection_drug"
    223 = "last_time_attacked_weapon"
    224 = "last_time_physical_abuse"
    225 = "last_time_sex_abuse"
    226 = "last_time_emotional_abuse"
    227 = "last_time_ongoing_abuse"
    228 = "last_time_afraid_abuse"
    229 = "last_time_distressed_past"
    230 = "last_time_suicidal"
    231 = "ever_dx_hiv"
    232 = "ever_dx_hcv"
    233 = "ever_dx_hepb"
    234 = "ever_dx_std"
    235 = "ever_dx_tb"
    236 = "ever_dx_covid"
    237 = "first_dx_hiv"
    238 = "first_dx_hcv"
    239 = "first_dx_hepb"
    240 = "first_dx_std"
    241 = "first_dx_tb"
    242 = "first_dx_covid"
}
foreach ($row in $ws2Updates.Keys) {
    $ws2.Range("K$row").Value = $ws2Updates[$row]
}
$ws2.Columns.Item(11).ColumnWidth = 31.166666666666668
